# Weekly price-update edit: a new week's record for Acelga at
# "Terminal Hortofrutícola Agro Chillán" is inserted at the top of the
# data block (row 105), pushing all existing records (rows 105-160) down
# by one row (new rows 106-161).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the first data row of this block (row 105),
# shifting rows 105:160 down to 106:161 (Excel copies formatting from the
# row above automatically, same as Excel UI "Insert Sheet Rows").
$ws.Rows("105:105").Insert()

# Populate the newly-inserted row 105 with the new record.
$ws.Range("A105").Value = 7
$ws.Range("B105").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C105").Value = "Ñuble"
$ws.Range("D105").Value = 44460
$ws.Range("E105").Value = 16
$ws.Range("F105").Value = 100112009
$ws.Range("G105").Value = "Acelga"
$ws.Range("H105").Value = "Sin especificar"
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 100
$ws.Range("K105").Value = 350
$ws.Range("L105").Value = 400
$ws.Range("M105").Value = 375
$ws.Range("N105").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O105").Value = "Provincia de Diguillín"
$ws.Range("P105").Value = 375
$ws.Range("Q105").Value = 1
$ws.Range("R105").Value = "Hortaliza"
